# "Hortaliza, Macroferia Regional de Talca - Pimiento" weekly update.
#
# The new weekly price report adds two rows of data (Zafiro rojo / Zafiro
# verde @ 2022-05-30) right after the existing header block's first data
# row (row 439), pushing every existing data row down by two. The sheet's
# used range grows from A1:R535 to A1:R537 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 439:440 - this shifts old rows 439..535 down to
# 441..537 and extends the sheet dimension automatically.
$ws.Rows("439:440").Insert()

# New row 439: Zafiro rojo, 2022-05-30
$ws.Range("A439").Value = 5
$ws.Range("B439").Value = "Macroferia Regional de Talca"
$ws.Range("C439").Value = "Maule"
$ws.Range("D439").Value = 44711
$ws.Range("E439").Value = 7
$ws.Range("F439").Value = 100112002
$ws.Range("G439").Value = "Pimiento"
$ws.Range("H439").Value = "Zafiro rojo"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 300
$ws.Range("K439").Value = 45000
$ws.Range("L439").Value = 45000
$ws.Range("M439").Value = 45000
$ws.Range("N439").Value = "$/caja 15 kilos"
$ws.Range("O439").Value = "Región de Arica y Parinacota"
$ws.Range("P439").Value = 3000
$ws.Range("Q439").Value = 15
$ws.Range("R439").Value = "Hortaliza"

# New row 440: Zafiro verde, 2022-05-30
$ws.Range("A440").Value = 5
$ws.Range("B440").Value = "Macroferia Regional de Talca"
$ws.Range("C440").Value = "Maule"
$ws.Range("D440").Value = 44711
$ws.Range("E440").Value = 7
$ws.Range("F440").Value = 100112002
$ws.Range("G440").Value = "Pimiento"
$ws.Range("H440").Value = "Zafiro verde"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 300
$ws.Range("K440").Value = 25000
$ws.Range("L440").Value = 25000
$ws.Range("M440").Value = 25000
$ws.Range("N440").Value = "$/caja 15 kilos"
$ws.Range("O440").Value = "Región de Arica y Parinacota"
$ws.Range("P440").Value = 1667
$ws.Range("Q440").Value = 15
$ws.Range("R440").Value = "Hortaliza"
